$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '29.035.93'
$ws.Cells.Item(2, 5).Value = '  +0.21%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.921.95'
$ws.Cells.Item(3, 5).Value = '  +0.95%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.39%  '

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '325.45'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.41%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +0.43%  '

# Row 7
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.4598'
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +0.06%  '

# Row 8
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '0.3816'
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.29%  '

# Row 9
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.07747'
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +0.01%  '

# Row 10
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.9789'
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +0.04%  '

# Row 11
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '22.75'
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +2.78%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.939.01'
$ws.Cells.Item(12, 5).Value = '  +2.78%  '

# Row 13
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '5.704'
$c.Style = "Normal"

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '6.963'
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -0.32%  '

# Row 15
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '0.07018'
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -0.32%  '

# Row 16
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '1.007'
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.50%  '

# Row 17
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '84.53'
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.81%  '

# Row 18
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '0.000009499'
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.36%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -0.17%  '

# Row 20
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.42%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '29.046.10'
$ws.Cells.Item(21, 5).Value = '  +0.21%  '

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '5.356'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.56%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +1.01%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '2.159.48'
$ws.Cells.Item(24, 5).Value = '  +1.59%  '

# Row 25
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '2.072'
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.14%  '

# Row 26
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '157.79'
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +1.14%  '

# Row 27
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '19.00'
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.71%  '

# Row 28
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '5.634'
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +0.30%  '

# Row 29
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '117.94'
$c.Style = "Normal"

# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.30%  '

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '0.09306'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.52%  '

# Row 32
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '0.8605'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +0.31%  '

# Row 33
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '5.096'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -0.23%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -0.76%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +0.33%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '1.160'
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +0.90%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Hedera'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '0.05690'
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -0.60%  '

# Row 38
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +0.36%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'MXToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '3.140'
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +15.52%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '0.02046'
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.01%  '

# Row 41
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '7.420'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -0.54%  '

# Row 42
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '0.5490'
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.63%  '

# Row 44
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '9.383'
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +1.00%  '

# Row 45
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.000002769'
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -2.21%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +4.47%  '

# Row 47
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '0.5190'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -0.24%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '11.24'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.39%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '0.06926'
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +1.50%  '

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '110.32'
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -1.33%  '

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '1.761'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -1.17%  '
